$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 439.11765
$ws.Range("J17").Value = 439.11765
$ws.Range("L17").Value = 1317.35295
$ws.Range("N17").Value = -1653.35295
$ws.Range("H43").Value = 841.6667
$ws.Range("I43").Value = 491
$ws.Range("J43").Value = 885.5
$ws.Range("K43").Value = 491
$ws.Range("L43").Value = 885.5
$ws.Range("M43").Value = -422
$ws.Range("N43").Value = -1023.5
$ws.Range("H112").Value = 1257.4706
$ws.Range("J112").Value = 1287.5
$ws.Range("L112").Value = 3862.5
$ws.Range("N112").Value = -6078.5
$ws.Range("H127").Value = 111111784
$ws.Range("I127").Value = 125000570
$ws.Range("J127").Value = 1500
$ws.Range("K127").Value = 375001710
$ws.Range("L127").Value = 4500
$ws.Range("M127").Value = -374996750
$ws.Range("N127").Value = -14420

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 2788.6
$ws.Range("I3").Value = 252.5
$ws.Range("K3").Value = 252.5
$ws.Range("M3").Value = -137.5
$ws.Range("H61").Value = 1247.9166
$ws.Range("I61").Value = 1076.1538
$ws.Range("J61").Value = 1450.909
$ws.Range("K61").Value = 1076.1538
$ws.Range("L61").Value = 1450.909
$ws.Range("M61").Value = -864.1538
$ws.Range("N61").Value = -1874.909
$ws.Range("H74").Value = 31505.084
$ws.Range("I74").Value = 40834.36
$ws.Range("J74").Value = 10302.182
$ws.Range("K74").Value = 40834.36
$ws.Range("L74").Value = 10302.182
$ws.Range("M74").Value = -39960.36
$ws.Range("N74").Value = -12050.182
$ws.Range("H77").Value = 31505.084
$ws.Range("I77").Value = 40834.36
$ws.Range("J77").Value = 10302.182
$ws.Range("K77").Value = 204171.8
$ws.Range("L77").Value = 51510.91
$ws.Range("M77").Value = -199803.8
$ws.Range("N77").Value = -60246.91
$ws.Range("H122").Value = 1673.2858
$ws.Range("I122").Value = 1473.7727
$ws.Range("J122").Value = 2010.9231
$ws.Range("K122").Value = 4421.3181
$ws.Range("L122").Value = 6032.7693
$ws.Range("M122").Value = -1971.3181
$ws.Range("N122").Value = -10932.7693
$ws.Range("H132").Value = 2361.375
$ws.Range("I132").Value = 1881.1666
$ws.Range("J132").Value = 3802
$ws.Range("K132").Value = 5643.4998
$ws.Range("L132").Value = 11406
$ws.Range("M132").Value = -3113.4998
$ws.Range("N132").Value = -16466
$ws.Range("H136").Value = 1247.9166
$ws.Range("I136").Value = 1076.1538
$ws.Range("J136").Value = 1450.909
$ws.Range("K136").Value = 3228.4614
$ws.Range("L136").Value = 4352.727000000001
$ws.Range("M136").Value = -678.4614000000001
$ws.Range("N136").Value = -9452.727000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2479.8538
$ws.Range("J105").Value = 2494.7334
$ws.Range("L105").Value = 2494.7334
$ws.Range("N105").Value = -5988.7334
$ws.Range("H134").Value = 5528.271
$ws.Range("I134").Value = 5734.4517
$ws.Range("J134").Value = 5152.294
$ws.Range("K134").Value = 17203.3551
$ws.Range("L134").Value = 15456.882
$ws.Range("M134").Value = -14668.3551
$ws.Range("N134").Value = -20526.882

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15153331
$ws.Range("I31").Value = 1553.5
$ws.Range("J31").Value = 21741060
$ws.Range("K31").Value = 1553.5
$ws.Range("L31").Value = 21741060
$ws.Range("M31").Value = -1258.5
$ws.Range("N31").Value = -21741650
$ws.Range("H34").Value = 15153331
$ws.Range("I34").Value = 1553.5
$ws.Range("J34").Value = 21741060
$ws.Range("K34").Value = 1553.5
$ws.Range("L34").Value = 21741060
$ws.Range("M34").Value = -1351.5
$ws.Range("N34").Value = -21741464
$ws.Range("H58").Value = 5315.826
$ws.Range("I58").Value = 6314.2104
$ws.Range("J58").Value = 573.5
$ws.Range("K58").Value = 6314.2104
$ws.Range("L58").Value = 573.5
$ws.Range("M58").Value = -6111.2104
$ws.Range("N58").Value = -979.5
$ws.Range("H112").Value = 29873.857
$ws.Range("J112").Value = 29873.857
$ws.Range("L112").Value = 29873.857
$ws.Range("N112").Value = -32827.857
$ws.Range("H132").Value = 3142
$ws.Range("I132").Value = 2804.25
$ws.Range("J132").Value = 3817.5
$ws.Range("K132").Value = 8412.75
$ws.Range("L132").Value = 11452.5
$ws.Range("M132").Value = -5882.75
$ws.Range("N132").Value = -16512.5
$ws.Range("H134").Value = 37933084
$ws.Range("I134").Value = 4547506.5
$ws.Range("J134").Value = 142859180
$ws.Range("K134").Value = 13642519.5
$ws.Range("L134").Value = 428577540
$ws.Range("M134").Value = -13639984.5
$ws.Range("N134").Value = -428582610
$ws.Range("H136").Value = 5315.826
$ws.Range("I136").Value = 6314.2104
$ws.Range("J136").Value = 573.5
$ws.Range("K136").Value = 18942.6312
$ws.Range("L136").Value = 1720.5
$ws.Range("M136").Value = -16392.6312
$ws.Range("N136").Value = -6820.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 366.66666
$ws.Range("I41").Value = 550
$ws.Range("J41").Value = 275
$ws.Range("K41").Value = 1650
$ws.Range("L41").Value = 825
$ws.Range("M41").Value = -1312
$ws.Range("N41").Value = -1501
$ws.Range("H48").Value = 1911.6666
$ws.Range("I48").Value = 1000
$ws.Range("K48").Value = 3000
$ws.Range("M48").Value = -2750
$ws.Range("H113").Value = 531.2646999999999
$ws.Range("I113").Value = 550.17645
$ws.Range("J113").Value = 512.35297
$ws.Range("K113").Value = 1650.52935
$ws.Range("L113").Value = 1537.05891
$ws.Range("M113").Value = 519.4706499999998
$ws.Range("N113").Value = -5877.05891
$ws.Range("H131").Value = 888.4262
$ws.Range("J131").Value = 1001.7059
$ws.Range("L131").Value = 3005.1177
$ws.Range("N131").Value = -13085.1177
$ws.Range("H134").Value = 6431.654
$ws.Range("I134").Value = 3744.5454
$ws.Range("J134").Value = 8402.200000000001
$ws.Range("K134").Value = 11233.6362
$ws.Range("L134").Value = 25206.6
$ws.Range("M134").Value = -6163.636200000001
$ws.Range("N134").Value = -35346.60000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1775
$ws.Range("I126").Value = 1800
$ws.Range("K126").Value = 5400
$ws.Range("M126").Value = -2930
$ws.Range("H132").Value = 6701.4287
$ws.Range("I132").Value = 7978
$ws.Range("J132").Value = 4999.3335
$ws.Range("K132").Value = 23934
$ws.Range("L132").Value = 14998.0005
$ws.Range("M132").Value = -21404
$ws.Range("N132").Value = -20058.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6721.6665
$ws.Range("I40").Value = 6512.5
$ws.Range("J40").Value = 7319.2856
$ws.Range("K40").Value = 6512.5
$ws.Range("L40").Value = 7319.2856
$ws.Range("M40").Value = -6376.5
$ws.Range("N40").Value = -7591.2856
$ws.Range("H93").Value = 71800.7
$ws.Range("I93").Value = 4334.3335
$ws.Range("J93").Value = 100714.86
$ws.Range("K93").Value = 4334.3335
$ws.Range("L93").Value = 100714.86
$ws.Range("M93").Value = -3086.3335
$ws.Range("N93").Value = -103210.86
$ws.Range("H122").Value = 3473
$ws.Range("I122").Value = 3488.3333
$ws.Range("J122").Value = 3450
$ws.Range("K122").Value = 10464.9999
$ws.Range("L122").Value = 10350
$ws.Range("M122").Value = -8014.999899999999
$ws.Range("N122").Value = -15250
$ws.Range("H132").Value = 2614.1096
$ws.Range("I132").Value = 2480.5813
$ws.Range("J132").Value = 2805.5
$ws.Range("K132").Value = 7441.743899999999
$ws.Range("L132").Value = 8416.5
$ws.Range("M132").Value = -4911.743899999999
$ws.Range("N132").Value = -13476.5
$ws.Range("H136").Value = 1013.5139
$ws.Range("I136").Value = 917.9259
$ws.Range("J136").Value = 1300.2778
$ws.Range("K136").Value = 2753.7777
$ws.Range("L136").Value = 3900.8334
$ws.Range("M136").Value = -203.7776999999996
$ws.Range("N136").Value = -9000.8334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 30000
$ws.Range("J112").Value = 30000
$ws.Range("L112").Value = 30000
$ws.Range("N112").Value = -32954
$ws.Range("H126").Value = 55557030
$ws.Range("I126").Value = 1120.9231
$ws.Range("K126").Value = 3362.7693
$ws.Range("M126").Value = -892.7692999999999
$ws.Range("H132").Value = 6438.1577
$ws.Range("I132").Value = 9130.9
$ws.Range("J132").Value = 3446.2222
$ws.Range("K132").Value = 27392.7
$ws.Range("L132").Value = 10338.6666
$ws.Range("M132").Value = -24862.7
$ws.Range("N132").Value = -15398.6666
$ws.Range("H136").Value = 2491.1082
$ws.Range("I136").Value = 2462.8147
$ws.Range("J136").Value = 2567.5
$ws.Range("K136").Value = 7388.4441
$ws.Range("L136").Value = 7702.5
$ws.Range("M136").Value = -4838.4441
$ws.Range("N136").Value = -12802.5
